$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.028.89"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.418.94"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'563.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").Value = "'142.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  -3.09%  "
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("D13").Value = "'25.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "'0.0000173"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "2.854.97"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "61.993.40"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "2.424.29"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "'11.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("D19").Value = "'321.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "'6.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "'66.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").Value = "'8.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.65%  "
$ws.Range("D26").Value = "'563.93"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("D27").Value = "2.533.89"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").Value = "'1.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "'4.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.03%  "
$ws.Range("D37").Value = "'5.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.62%  "
$ws.Range("D38").Value = "'153.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.41%  "
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "'18.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.00%  "
$ws.Range("D41").Value = "'1.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.21%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'149.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").Value = "'19.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.30%  "
$ws.Range("D48").Value = "'0.593"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("D49").Value = "'0.0921"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  -0.47%  "
